# Update the "量的変数" (quantitative variables) sheet:
#   - D2 (max/Pressure row) changes from 5 to 4
#   - the saved cursor selection moves from C3 to G7
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("量的変数")
$ws.Activate()

$ws.Range("D2").Value = 4

$ws.Range("G7").Select()
